$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new transaction row (row 42) below the existing ledger data,
# matching the formatting of the existing rows.

# A42: date (same number format as the date column above it)
$ws.Range("A42").Value = 44386
$ws.Range("A42").NumberFormat = $ws.Range("A41").NumberFormat

# B42: amount added
$ws.Range("B42").Value = 15290

# C42: wallet type ("Connectivity"), centered like the rest of column C
$ws.Range("C42").Value = "Connectivity"
$ws.Range("C42").HorizontalAlignment = -4108

# D42: remarks - copy formatting from a similarly-highlighted PLI commission
# remark cell so it picks up the existing highlight style instead of a new one
$ws.Range("D42").Value = "PLI Commission for June 2021 Vendor ledger"
$ws.Range("D10").Copy()
$ws.Range("D42").PasteSpecial(-4122)

$ws.Range("D42").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
